$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F88").Value = 75
$ws.Range("G88").Value = 8634
$ws.Range("F92").Value = 63
$ws.Range("G92").Value = 7052.22
$ws.Range("F97").Value = 360
$ws.Range("G97").Value = 16826.4
$ws.Range("F109").Value = 113
$ws.Range("G109").Value = 14178.11
$ws.Range("B114").Value = 245881.84
$ws.Range("B163").Value = 57552
$ws.Range("E163").Value = 136.86
$ws.Range("F163").Value = -5
$ws.Range("G163").Value = -603.45
$ws.Range("B164").Value = 64329
$ws.Range("E164").Value = 128.32
$ws.Range("F164").Value = 3
$ws.Range("G164").Value = 362.07
$ws.Range("F174").Value = 24
$ws.Range("G174").Value = 1800.48
$ws.Range("B178").Value = 3385.38
$ws.Range("F195").Value = 117
$ws.Range("G195").Value = 10164.96
$ws.Range("F196").Value = 91
$ws.Range("G196").Value = 8034.39
$ws.Range("F198").Value = 20
$ws.Range("G198").Value = 1502.8
$ws.Range("B200").Value = 44486.82
$ws.Range("F220").Value = 14
$ws.Range("G220").Value = 6324.36
$ws.Range("B222").Value = 45644.01
$ws.Range("F246").Value = 23
$ws.Range("G246").Value = 765.9
$ws.Range("F248").Value = 46
$ws.Range("G248").Value = 3943.12
$ws.Range("F250").Value = 21
$ws.Range("G250").Value = 2405.76
$ws.Range("F259").Value = 2
$ws.Range("G259").Value = 271.08
$ws.Range("B274").Value = 74005.32000000001
$ws.Range("B294").Value = 57802
$ws.Range("E294").Value = 162.71
$ws.Range("F294").Value = -79
$ws.Range("G294").Value = -11334.92
$ws.Range("B295").Value = 63571
$ws.Range("E295").Value = 152.53
$ws.Range("F295").Value = 2
$ws.Range("G295").Value = 286.96
$ws.Range("B296").Value = 63531
$ws.Range("F296").Value = 80
$ws.Range("G296").Value = 11478.4
$ws.Range("B299").Value = 55356
$ws.Range("E299").Value = 54.04
$ws.Range("F299").Value = -158
$ws.Range("G299").Value = -7527.12
$ws.Range("B300").Value = 63510
$ws.Range("E300").Value = 50.66
$ws.Range("F300").Value = 135
$ws.Range("G300").Value = 6431.4
$ws.Range("F310").Value = 51
$ws.Range("G310").Value = 2658.63
$ws.Range("F328").Value = 494
$ws.Range("G328").Value = 10388.82
$ws.Range("F333").Value = 682
$ws.Range("G333").Value = 116847.06
$ws.Range("F334").Value = 351
$ws.Range("G334").Value = 53060.67
$ws.Range("B339").Value = 282546.59
$ws.Range("F355").Value = 80
$ws.Range("G355").Value = 2579.2
$ws.Range("B361").Value = 12103.5
$ws.Range("F363").Value = 28
$ws.Range("G363").Value = 595
$ws.Range("F368").Value = 157
$ws.Range("G368").Value = 26554.98
$ws.Range("F372").Value = 14
$ws.Range("G372").Value = 755.3
$ws.Range("F384").Value = 64
$ws.Range("G384").Value = 4500.48
$ws.Range("F385").Value = 77
$ws.Range("G385").Value = 5494.72
$ws.Range("F386").Value = 17
$ws.Range("G386").Value = 1762.73
$ws.Range("F390").Value = 282
$ws.Range("G390").Value = 40788.48
$ws.Range("F393").Value = 108
$ws.Range("G393").Value = 22776.12
$ws.Range("B395").Value = 235564.58
$ws.Range("F398").Value = 87
$ws.Range("G398").Value = 1836.57
$ws.Range("B402").Value = 5033.04
$ws.Range("B420").Value = 58047
$ws.Range("D420").Value = 105.54
$ws.Range("E420").Value = 126.1
$ws.Range("F420").Value = 42
$ws.Range("G420").Value = 4432.68
$ws.Range("B421").Value = 47097
$ws.Range("D421").Value = 112.28
$ws.Range("E421").Value = 134.16
$ws.Range("F421").Value = 15
$ws.Range("G421").Value = 1684.2
$ws.Range("F426").Value = 65
$ws.Range("G426").Value = 6279
$ws.Range("B430").Value = 42532.58
$ws.Range("F446").Value = 138
$ws.Range("G446").Value = 8474.58
$ws.Range("B448").Value = 37983.71
$ws.Range("F470").Value = 96
$ws.Range("G470").Value = 1229.76
$ws.Range("F471").Value = 62
$ws.Range("G471").Value = 1630.6
$ws.Range("F479").Value = 135
$ws.Range("G479").Value = 2189.7
$ws.Range("F484").Value = 453
$ws.Range("G484").Value = 2939.97
$ws.Range("B485").Value = 64925
$ws.Range("E485").Value = 13.97
$ws.Range("F485").Value = 157
$ws.Range("G485").Value = 2064.55
$ws.Range("B486").Value = 45709
$ws.Range("E486").Value = 15.69
$ws.Range("F486").Value = -300
$ws.Range("G486").Value = -3945
$ws.Range("F488").Value = 67
$ws.Range("G488").Value = 1762.1
$ws.Range("F489").Value = 32
$ws.Range("G489").Value = 525.76
$ws.Range("B492").Value = -11807.78
$ws.Range("F497").Value = 20
$ws.Range("G497").Value = 1001.6
$ws.Range("B508").Value = 10132.74
$ws.Range("F551").Value = 117
$ws.Range("G551").Value = 1508.13
$ws.Range("F552").Value = 48
$ws.Range("G552").Value = 952.8
$ws.Range("B557").Value = 6699.57
$ws.Range("B601").Value = 60022
$ws.Range("E601").Value = 37.22
$ws.Range("F601").Value = -113
$ws.Range("G601").Value = -3709.79
$ws.Range("B602").Value = 64830
$ws.Range("E602").Value = 34.9
$ws.Range("F602").Value = 111
$ws.Range("G602").Value = 3644.13
$ws.Range("F621").Value = 26
$ws.Range("G621").Value = 3394.3
$ws.Range("B626").Value = 19577.48
$ws.Range("F703").Value = 13
$ws.Range("G703").Value = 1060.28
$ws.Range("F707").Value = 20
$ws.Range("G707").Value = 1631.2
$ws.Range("F712").Value = 41
$ws.Range("G712").Value = 890.52
$ws.Range("F715").Value = 77
$ws.Range("G715").Value = 5356.12
$ws.Range("F718").Value = 263
$ws.Range("G718").Value = 35507.63
$ws.Range("F720").Value = 185
$ws.Range("G720").Value = 22331.35
$ws.Range("B721").Value = 102765.87
$ws.Range("F725").Value = 24
$ws.Range("G725").Value = 3933.36
$ws.Range("F740").Value = 76
$ws.Range("G740").Value = 2512.56
$ws.Range("F742").Value = 135
$ws.Range("G742").Value = 5803.65
$ws.Range("F743").Value = 131
$ws.Range("G743").Value = 6535.59
$ws.Range("F746").Value = 92
$ws.Range("G746").Value = 22245.6
$ws.Range("F747").Value = 105
$ws.Range("G747").Value = 5983.95
$ws.Range("B748").Value = 66910.77
$ws.Range("F773").Value = 2916
$ws.Range("G773").Value = 475628.76
$ws.Range("F774").Value = 33
$ws.Range("G774").Value = 5813.61
$ws.Range("F775").Value = 570
$ws.Range("G775").Value = 161235.9
$ws.Range("F779").Value = 114
$ws.Range("G779").Value = 14658.12
$ws.Range("B780").Value = 724076.11
$ws.Range("F783").Value = 72
$ws.Range("G783").Value = 10512.72
$ws.Range("F790").Value = 192
$ws.Range("G790").Value = 30558.72
$ws.Range("F791").Value = 8
$ws.Range("G791").Value = 1284.08
$ws.Range("B797").Value = 70560.02
$ws.Range("B798").Value = 2722486.86
$ws.Range("B799").Value = 2722486.86
